# Auto-generated Excel COM-interop script
# Applies the 2024-01-21 daily data-refresh update to the violent-crime workbook.
# For each neighborhood/citywide sheet, updates the affected 2024 ("K" column, and a few
# other historical-year cells whose year totals were revised) crime-count cells.

$wb = $excel.ActiveWorkbook

# Sheet 1: Citywide Totals
$ws = $wb.Worksheets.Item(1)
$ws.Range("K2").Value = 346
$ws.Range("K3").Value = 309
$ws.Range("G4").Value = 1478
$ws.Range("J4").Value = 1772
$ws.Range("K4").Value = 65
$ws.Range("J6").Value = 11051
$ws.Range("K6").Value = 424
$ws.Range("G7").Value = 24702
$ws.Range("J7").Value = 29220
$ws.Range("K7").Value = 1162

# Sheet 2: By Neighborhood
$ws = $wb.Worksheets.Item(2)
$ws.Range("J5").Value = 87
$ws.Range("I7").Value = 822
$ws.Range("J7").Value = 823
$ws.Range("K7").Value = 33
$ws.Range("K8").Value = 79
$ws.Range("K11").Value = 34
$ws.Range("K16").Value = 4
$ws.Range("K18").Value = 6
$ws.Range("I19").Value = 732
$ws.Range("J19").Value = 858
$ws.Range("K19").Value = 22
$ws.Range("K20").Value = 35
$ws.Range("K23").Value = 10
$ws.Range("K29").Value = 57
$ws.Range("J33").Value = 1320
$ws.Range("K33").Value = 52
$ws.Range("K34").Value = 8
$ws.Range("K37").Value = 31
$ws.Range("K48").Value = 14
$ws.Range("K50").Value = 3
$ws.Range("K51").Value = 21
$ws.Range("K53").Value = 13
$ws.Range("K54").Value = 18
$ws.Range("K55").Value = 13
$ws.Range("K60").Value = 9
$ws.Range("G63").Value = 278
$ws.Range("J63").Value = 85
$ws.Range("K65").Value = 32
$ws.Range("K67").Value = 44
$ws.Range("K79").Value = 25
$ws.Range("K83").Value = 19
$ws.Range("K85").Value = 59
$ws.Range("J88").Value = 314
$ws.Range("K88").Value = 17
$ws.Range("K89").Value = 17
$ws.Range("K90").Value = 8
$ws.Range("K92").Value = 5
$ws.Range("J94").Value = 332
$ws.Range("K94").Value = 10
$ws.Range("K96").Value = 14
$ws.Range("K97").Value = 10
$ws.Range("J98").Value = 217
$ws.Range("K99").Value = 25
$ws.Range("G101").Value = 24702
$ws.Range("J101").Value = 29220
$ws.Range("K101").Value = 1162

# Sheet 4: West Ridge
$ws = $wb.Worksheets.Item(4)
$ws.Range("K2").Value = 7
$ws.Range("K7").Value = 14

# Sheet 5: Auburn Gresham
$ws = $wb.Worksheets.Item(5)
$ws.Range("I5").Value = 35
$ws.Range("J5").Value = 22
$ws.Range("K6").Value = 8
$ws.Range("I7").Value = 822
$ws.Range("J7").Value = 823
$ws.Range("K7").Value = 33

# Sheet 6: Belmont Cragin
$ws = $wb.Worksheets.Item(6)
$ws.Range("K3").Value = 5
$ws.Range("K7").Value = 34

# Sheet 7: Uptown
$ws = $wb.Worksheets.Item(7)
$ws.Range("K2").Value = 5
$ws.Range("K4").Value = 3
$ws.Range("K7").Value = 17

# Sheet 8: South Shore
$ws = $wb.Worksheets.Item(8)
$ws.Range("K2").Value = 23
$ws.Range("K3").Value = 18
$ws.Range("K7").Value = 59

# Sheet 11: Logan Square
$ws = $wb.Worksheets.Item(11)
$ws.Range("K6").Value = 6
$ws.Range("K7").Value = 13

# Sheet 12: Austin
$ws = $wb.Worksheets.Item(12)
$ws.Range("K3").Value = 26
$ws.Range("K7").Value = 79

# Sheet 13: South Chicago
$ws = $wb.Worksheets.Item(13)
$ws.Range("K3").Value = 5
$ws.Range("K7").Value = 19

# Sheet 14: Garfield Park
$ws = $wb.Worksheets.Item(14)
$ws.Range("K3").Value = 14
$ws.Range("J6").Value = 475
$ws.Range("J7").Value = 1320
$ws.Range("K7").Value = 52

# Sheet 16: Grand Crossing
$ws = $wb.Worksheets.Item(16)
$ws.Range("K2").Value = 7
$ws.Range("K6").Value = 10
$ws.Range("K7").Value = 31

# Sheet 17: New City
$ws = $wb.Worksheets.Item(17)
$ws.Range("K2").Value = 8
$ws.Range("K3").Value = 7
$ws.Range("K7").Value = 32

# Sheet 18: Woodlawn
$ws = $wb.Worksheets.Item(18)
$ws.Range("K2").Value = 10
$ws.Range("K6").Value = 7
$ws.Range("K7").Value = 25

# Sheet 21: North Lawndale
$ws = $wb.Worksheets.Item(21)
$ws.Range("K3").Value = 15
$ws.Range("K7").Value = 44

# Sheet 24: Loop
$ws = $wb.Worksheets.Item(24)
$ws.Range("K3").Value = 6
$ws.Range("K7").Value = 18

# Sheet 25: Englewood
$ws = $wb.Worksheets.Item(25)
$ws.Range("K2").Value = 19
$ws.Range("K6").Value = 20
$ws.Range("K7").Value = 57

# Sheet 26: Lake View
$ws = $wb.Worksheets.Item(26)
$ws.Range("K3").Value = 4
$ws.Range("K4").Value = 4
$ws.Range("K7").Value = 14

# Sheet 27: Chatham
$ws = $wb.Worksheets.Item(27)
$ws.Range("K3").Value = 8
$ws.Range("I5").Value = 24
$ws.Range("J5").Value = 31
$ws.Range("K6").Value = 6
$ws.Range("I7").Value = 732
$ws.Range("J7").Value = 858
$ws.Range("K7").Value = 22

# Sheet 36: Lower West Side
$ws = $wb.Worksheets.Item(36)
$ws.Range("K6").Value = 5
$ws.Range("K7").Value = 13

# Sheet 39: Douglas
$ws = $wb.Worksheets.Item(39)
$ws.Range("K4").Value = 1
$ws.Range("K7").Value = 10

# Sheet 42: Roseland
$ws = $wb.Worksheets.Item(42)
$ws.Range("K2").Value = 7
$ws.Range("K6").Value = 8
$ws.Range("K7").Value = 25

# Sheet 44: Chicago Lawn
$ws = $wb.Worksheets.Item(44)
$ws.Range("K3").Value = 9
$ws.Range("K6").Value = 18
$ws.Range("K7").Value = 35

# Sheet 45: Calumet Heights
$ws = $wb.Worksheets.Item(45)
$ws.Range("K3").Value = 2
$ws.Range("K7").Value = 6

# Sheet 50: Garfield Ridge
$ws = $wb.Worksheets.Item(50)
$ws.Range("K4").Value = 1
$ws.Range("K7").Value = 8

# Sheet 51: West Loop
$ws = $wb.Worksheets.Item(51)
$ws.Range("K2").Value = 3
$ws.Range("K3").Value = 1
$ws.Range("J4").Value = 28
$ws.Range("K6").Value = 6
$ws.Range("J7").Value = 332
$ws.Range("K7").Value = 10

# Sheet 55: Wicker Park
$ws = $wb.Worksheets.Item(55)
$ws.Range("J4").Value = 11
$ws.Range("J7").Value = 217

# Sheet 56: Lincoln Square
$ws = $wb.Worksheets.Item(56)
$ws.Range("K6").Value = 3
$ws.Range("K7").Value = 3

# Sheet 65: West Town
$ws = $wb.Worksheets.Item(65)
$ws.Range("K6").Value = 4
$ws.Range("K7").Value = 10

# Sheet 66: West Elsdon
$ws = $wb.Worksheets.Item(66)
$ws.Range("K2").Value = 1
$ws.Range("K7").Value = 5

# Sheet 68: United Center
$ws = $wb.Worksheets.Item(68)
$ws.Range("K2").Value = 5
$ws.Range("J4").Value = 10
$ws.Range("J7").Value = 314
$ws.Range("K7").Value = 17

# Sheet 70: Armour Square
$ws = $wb.Worksheets.Item(70)
$ws.Range("J4").Value = 2
$ws.Range("J7").Value = 87

# Sheet 74: Washington Heights
$ws = $wb.Worksheets.Item(74)
$ws.Range("K3").Value = 4
$ws.Range("K7").Value = 8

# Sheet 75: Little Italy, UIC
$ws = $wb.Worksheets.Item(75)
$ws.Range("K2").Value = 8
$ws.Range("K7").Value = 21

# Sheet 78: Morgan Park
$ws = $wb.Worksheets.Item(78)
$ws.Range("K6").Value = 1
$ws.Range("K7").Value = 9

# Sheet 94: Bucktown
$ws = $wb.Worksheets.Item(94)
$ws.Range("K6").Value = 3
$ws.Range("K7").Value = 4
